$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "'1231"
$ws.Range("D2").Value = "'23123"
$ws.Range("F2").Value = "'12312"
$ws.Range("G2").Value = "'3123"
$ws.Range("H2").Value = "'1231"
